$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, date range) ---
$ws.Range("A8").Value = "Volume 30   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/12/2023  Through  6/18/2023"

# --- Type-changing cells: copy style/type from a stable same-shaped cell first ---
$ws.Range("C23").Copy($ws.Range("C14"))
$ws.Range("C16").Copy($ws.Range("D15"))
$ws.Range("E16").Copy($ws.Range("E15"))
$ws.Range("C23").Copy($ws.Range("D22"))
$ws.Range("E23").Copy($ws.Range("E22"))
$ws.Range("C16").Copy($ws.Range("D26"))
$ws.Range("E16").Copy($ws.Range("E26"))
$ws.Range("C23").Copy($ws.Range("C27"))
$ws.Range("C16").Copy($ws.Range("D28"))
$ws.Range("E16").Copy($ws.Range("E28"))
$ws.Range("C16").Copy($ws.Range("D29"))
$ws.Range("E16").Copy($ws.Range("E29"))

# --- Now set the actual values for the type-changed cells ---
$ws.Range("C14").Value = "0"
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -50
$ws.Range("D22").Value = "0"
$ws.Range("E22").Value = "***.*"
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -50
$ws.Range("C27").Value = "0"
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 0

# --- Plain value updates (style/type unchanged) ---
$ws.Range("N14").Value = -55.555555555555
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = -16.666666666666
$ws.Range("I15").Value = 19
$ws.Range("J15").Value = 16
$ws.Range("K15").Value = 18.75
$ws.Range("L15").Value = 35.714285714285
$ws.Range("M15").Value = 11.764705882352
$ws.Range("N15").Value = -55.813953488372
$ws.Range("C16").Value = 7
$ws.Range("E16").Value = -12.5
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 31
$ws.Range("H16").Value = -41.935483870967
$ws.Range("I16").Value = 121
$ws.Range("J16").Value = 167
$ws.Range("K16").Value = -27.54491017964
$ws.Range("L16").Value = 27.368421052631
$ws.Range("M16").Value = -35.978835978836
$ws.Range("N16").Value = -89.128481581311
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = -35.294117647058
$ws.Range("F17").Value = 48
$ws.Range("G17").Value = 59
$ws.Range("H17").Value = -18.64406779661
$ws.Range("I17").Value = 304
$ws.Range("J17").Value = 283
$ws.Range("K17").Value = 7.420494699646
$ws.Range("L17").Value = 15.151515151515
$ws.Range("M17").Value = 49.019607843137
$ws.Range("N17").Value = -44.727272727272
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 36.363636363636
$ws.Range("I18").Value = 98
$ws.Range("J18").Value = 94
$ws.Range("K18").Value = 4.255319148936
$ws.Range("L18").Value = 60.655737704918
$ws.Range("M18").Value = -44.632768361581
$ws.Range("N18").Value = -89.705882352941
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 28.571428571428
$ws.Range("F19").Value = 53
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = -1.851851851851
$ws.Range("I19").Value = 323
$ws.Range("J19").Value = 331
$ws.Range("K19").Value = -2.416918429003
$ws.Range("L19").Value = 78.453038674033
$ws.Range("M19").Value = 28.685258964143
$ws.Range("N19").Value = -6.376811594202
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -50
$ws.Range("I20").Value = 119
$ws.Range("J20").Value = 135
$ws.Range("K20").Value = -11.851851851851
$ws.Range("L20").Value = 19
$ws.Range("M20").Value = -9.848484848484
$ws.Range("N20").Value = -87.943262411347
$ws.Range("C21").Value = 42
$ws.Range("E21").Value = -17.647058823529
$ws.Range("F21").Value = 158
$ws.Range("G21").Value = 189
$ws.Range("H21").Value = -16.402116402116
$ws.Range("I21").Value = 992
$ws.Range("J21").Value = 1035
$ws.Range("K21").Value = -4.15458937198
$ws.Range("L21").Value = 37.777777777777
$ws.Range("M21").Value = 1.535312180143
$ws.Range("N21").Value = -75.249500998004
$ws.Range("G22").Value = 1
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 92.857142857142
$ws.Range("F24").Value = 75
$ws.Range("G24").Value = 75
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 566
$ws.Range("J24").Value = 436
$ws.Range("K24").Value = 29.816513761467
$ws.Range("L24").Value = 35.08353221957
$ws.Range("M24").Value = 29.223744292237
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = -39.130434782608
$ws.Range("F25").Value = 57
$ws.Range("G25").Value = 72
$ws.Range("H25").Value = -20.833333333333
$ws.Range("I25").Value = 350
$ws.Range("J25").Value = 334
$ws.Range("K25").Value = 4.790419161676
$ws.Range("L25").Value = 22.807017543859
$ws.Range("M25").Value = -16.067146282973
$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 7
$ws.Range("H26").Value = -28.571428571428
$ws.Range("I26").Value = 23
$ws.Range("J26").Value = 23
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 4.545454545454
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 6
$ws.Range("J27").Value = 39
$ws.Range("K27").Value = -10.25641025641
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -40
$ws.Range("I28").Value = 23
$ws.Range("J28").Value = 21
$ws.Range("K28").Value = 9.523809523809
$ws.Range("L28").Value = -8
$ws.Range("M28").Value = -8
$ws.Range("N28").Value = -74.157303370786
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -25
$ws.Range("I29").Value = 19
$ws.Range("J29").Value = 16
$ws.Range("K29").Value = 18.75
$ws.Range("L29").Value = -17.391304347826
$ws.Range("M29").Value = -9.523809523809
$ws.Range("N29").Value = -77.647058823529

